# Swap the header labels in B1 and C1 ("SP500 weight" <-> "Tbill weight")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b1 = $ws.Range("B1").Value()
$c1 = $ws.Range("C1").Value()

$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# Move the active selection to C2, matching the saved view state
$ws.Range("C2").Select()
